$d = $word.ActiveDocument

# --- Change 1: merge the grammar-check-split runs in the "Course registration"
#     paragraph back into a single run (drop the proofErr gramStart/gramEnd
#     markers that wrapped "in" and "the courses"). The visible text is
#     unchanged; a plain Find/Replace of the exact same text causes Word to
#     rewrite the whole matched span as a single run.
$oldCourseText = "The student can access his/her academic profile in the university" + [char]8217 + "s platform. The student can see what the courses he/she can register for with the specified credit hours. Then the student can add the courses wants to register for and submit the request. The request is then saved in a temporary file for the administrator to revise it and confirm it. After that the student will be notified with the acceptance or rejection message. If the student request is "
$d.Content.Find.Execute($oldCourseText, $true, $false, $false, $false, $false, $true, 1, $false, $oldCourseText, 2) | Out-Null

# --- Change 2: "Managing student records" + ":" runs -> merge into a single
#     run with text "Managing student records:".
$d.Content.Find.Execute("Managing student records:", $true, $false, $false, $false, $false, $true, 1, $false, "Managing student records:", 2) | Out-Null

# --- Change 3: split the "Managing student records:" heading paragraph from
#     the following body paragraph by inserting a new, empty paragraph
#     between them.
$heading = $d.Paragraphs(18)
$headingRange = $heading.Range
$headingRange.Collapse(0)
$headingRange.InsertAfter([char]13)

# --- Change 4: rewrite the tail of the "drop a course" paragraph describing
#     the automatic system check instead of the manual administrator review.
#     The unchanged lead-in (" student can submit a request ") stays in its
#     own run, and the replaced tail becomes a separate run with identical
#     formatting. A transient Bold toggle forces Word to keep the two runs
#     split instead of re-merging them.
$oldTail = "and the request will be added to a temporary file that contains all the requests. Then the administrator can check the requests and approve it then the dropped course will be deleted from the student registered courses in the student record and the platform will be changed accordingly. The student will be notified."
$newTail = "and the request will be checked automatically by the system to see if it goes below the limit (limit is 2 courses) or not, and then update it."
$tailFound = $d.Content
$tailFound.Find.Execute($oldTail) | Out-Null
$tailRange = $d.Range($tailFound.Start, $tailFound.End)
$tailRange.Bold = 1
$tailRange.Text = $newTail
$tailRange.Bold = 0
